# Daily attendance processing - 2025-11-29 13:33:31
# Move the leading "System" (or "system") token in the "Recorded By" column
# to the end of the comma-separated list, for every data row where it is
# currently the first token.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -ne "") {
        $parts = $value -split ", "
        if ($parts.Count -gt 1 -and $parts[0] -eq "System") {
            $rest = $parts[1..($parts.Count - 1)]
            $newValue = ($rest + $parts[0]) -join ", "
            $cell.Value2 = $newValue
        }
    }
}
